# PropertyTypes.xlsx - drop the "TimeSpan of" demo row from the
# AllPropertyTypes worksheet (TimeSpan support removed, per commit message).
#
# Deleting the whole row shifts every row below it up by one, which is
# exactly what the target workbook shows: the "Enum of" row (old row 9)
# becomes row 8, ... the trailing "Anything / = / Doesn't Matter" row (old
# row 20) becomes row 19, and the sheet's used range shrinks from
# A1:D20 to A1:D19. Excel/the xlsx writer also drops the now-unreferenced
# "TimeSpan of" shared string and renumbers the remaining ones
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 holds "TimeSpan of" (C8) / the TimeSpan value (D8) - remove it
# entirely (not just clear its contents) so everything below shifts up.
$ws.Range("A8").EntireRow.Delete()

# Row 1 ("Specification" / "All simple property types") no longer needs
# its explicit 30pt height once re-laid-out; let Excel recompute it.
$ws.Rows.Item(1).AutoFit()

# Leave the selection where the author's save left it.
$ws.Range("D20").Select() | Out-Null
